$d = $word.ActiveDocument

$startTag = "<w:docDefaults>"
$endTag = "</w:docDefaults>"

$xml = $d.WordOpenXML

$startIdx = $xml.IndexOf($startTag)
$endIdx = $xml.IndexOf($endTag) + $endTag.Length

$newDocDefaults = "<w:docDefaults><w:rPrDefault><w:rPr><w:rFonts w:ascii=`"Arial`" w:cs=`"Arial`" w:eastAsia=`"Arial`" w:hAnsi=`"Arial`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"22`"/><w:lang w:val=`"fr`"/></w:rPr></w:rPrDefault><w:pPrDefault><w:pPr><w:spacing w:line=`"276`" w:lineRule=`"auto`"/></w:pPr></w:pPrDefault></w:docDefaults>"

$newXml = $xml.Substring(0, $startIdx) + $newDocDefaults + $xml.Substring($endIdx)

$d.WordOpenXML = $newXml
